# Edit: insert a new record row at row 296 for "Vega Monumental Concepción - Acelga",
# shifting the existing rows 296-374 down to 297-375, and fill the new row 296
# with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at 296; this shifts rows 296:374 down to 297:375
# and copies formatting from the row above (row 295), matching the existing
# date-formatted style used throughout column D.
$ws.Rows.Item(296).Insert()

# Fill the newly inserted row 296 with the new record's data.
$ws.Cells.Item(296, 1).Value = 11
$ws.Cells.Item(296, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(296, 3).Value = "Bíobío"
$ws.Cells.Item(296, 4).Value = 44985
$ws.Cells.Item(296, 5).Value = 8
$ws.Cells.Item(296, 6).Value = 100112009
$ws.Cells.Item(296, 7).Value = "Acelga"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Primera"
$ws.Cells.Item(296, 10).Value = 350
$ws.Cells.Item(296, 11).Value = 550
$ws.Cells.Item(296, 12).Value = 600
$ws.Cells.Item(296, 13).Value = 579
$ws.Cells.Item(296, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(296, 15).Value = "Región de Ñuble"
$ws.Cells.Item(296, 16).Value = 579
$ws.Cells.Item(296, 17).Value = 1
$ws.Cells.Item(296, 18).Value = "Hortaliza"
